# Updated cryptos list data (Price / Volume(1h) columns) to match the
# upstream GitHub Actions refresh described by the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    # Assigning a numeric-looking string via .Value lets Excel coerce it to a
    # Double (dropping trailing zeros / exact formatting), so force literal
    # text entry via a leading apostrophe, then restore the original cell
    # style (the apostrophe trick otherwise flips in a "quote prefix" style).
    $origStyle = $range.Style
    $range.Formula = "'" + $text
    $range.Style = $origStyle
}

$ws.Range("D2").Value = "26.851.17"
$ws.Range("E2").Value = "  -0.90%  "
$ws.Range("D3").Value = "1.868.32"
$ws.Range("E3").Value = "  +0.09%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("E6").Value = "  +0.11%  "
Set-TextValue $ws.Range("D7") "0.5080"
$ws.Range("E7").Value = "  -1.30%  "
Set-TextValue $ws.Range("D8") "0.3657"
$ws.Range("E8").Value = "  -2.71%  "
Set-TextValue $ws.Range("D9") "0.07185"
$ws.Range("E9").Value = "  +0.40%  "
Set-TextValue $ws.Range("D10") "0.8923"
$ws.Range("E10").Value = "  +0.35%  "
Set-TextValue $ws.Range("D11") "20.68"
Set-TextValue $ws.Range("D12") "0.07524"
$ws.Range("E12").Value = "  -1.02%  "
$ws.Range("D13").Value = "1.876.64"
$ws.Range("E13").Value = "  -1.51%  "
Set-TextValue $ws.Range("D14") "94.72"
$ws.Range("E14").Value = "  +5.73%  "
$ws.Range("E15").Value = "  -1.42%  "
$ws.Range("E16").Value = "  +0.15%  "
Set-TextValue $ws.Range("D17") "0.000008489"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("E18").Value = "  +0.58%  "
Set-TextValue $ws.Range("D19") "1.001"
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("D20").Value = "26.899.25"
$ws.Range("E20").Value = "  -0.84%  "
Set-TextValue $ws.Range("D21") "5.008"
$ws.Range("E21").Value = "  -0.39%  "
$ws.Range("D22").Value = "2.117.04"
$ws.Range("E22").Value = "  +0.34%  "
$ws.Range("E23").Value = "  -1.42%  "
Set-TextValue $ws.Range("D24") "6.380"
$ws.Range("E24").Value = "  -1.21%  "
Set-TextValue $ws.Range("D25") "148.05"
$ws.Range("E25").Value = "  +0.34%  "
$ws.Range("E26").Value = "  -3.06%  "
Set-TextValue $ws.Range("D28") "2.094"
$ws.Range("E28").Value = "  -0.69%  "
$ws.Range("E29").Value = "  +0.65%  "
Set-TextValue $ws.Range("D30") "4.697"
$ws.Range("E30").Value = "  +0.77%  "
Set-TextValue $ws.Range("D31") "4.729"
$ws.Range("E31").Value = "  +0.65%  "
Set-TextValue $ws.Range("D32") "0.09143"
$ws.Range("E32").Value = "  +0.11%  "
Set-TextValue $ws.Range("D33") "0.05066"
$ws.Range("E33").Value = "  -1.08%  "
Set-TextValue $ws.Range("D34") "0.7470"
$ws.Range("E34").Value = "  +2.80%  "
Set-TextValue $ws.Range("D35") "2.981"
$ws.Range("E35").Value = "  -2.77%  "
Set-TextValue $ws.Range("D36") "1.154"
$ws.Range("E36").Value = "  -0.25%  "
Set-TextValue $ws.Range("D37") "3.229"
$ws.Range("E37").Value = "  +5.81%  "
Set-TextValue $ws.Range("D38") "2.528"
$ws.Range("E38").Value = "  +0.75%  "
Set-TextValue $ws.Range("D39") "0.5591"
$ws.Range("E39").Value = "  +4.70%  "
Set-TextValue $ws.Range("D40") "0.01992"
$ws.Range("E40").Value = "  -2.29%  "
$ws.Range("E41").Value = "  +0.12%  "
Set-TextValue $ws.Range("D42") "6.607"
$ws.Range("E42").Value = "  +0.72%  "
Set-TextValue $ws.Range("D43") "115.36"
$ws.Range("E43").Value = "  -0.27%  "
Set-TextValue $ws.Range("D44") "8.588"
$ws.Range("E44").Value = "  +3.44%  "
$ws.Range("E45").Value = "  +0.57%  "
Set-TextValue $ws.Range("D46") "0.4733"
$ws.Range("E46").Value = "  +2.07%  "
Set-TextValue $ws.Range("D47") "1.000"
$ws.Range("E47").Value = "  +0.13%  "
Set-TextValue $ws.Range("D48") "10.14"
$ws.Range("E48").Value = "  +1.62%  "
Set-TextValue $ws.Range("D49") "1.565"
$ws.Range("E49").Value = "  -0.41%  "
Set-TextValue $ws.Range("D50") "36.90"
$ws.Range("E50").Value = "  +1.02%  "
Set-TextValue $ws.Range("D51") "63.01"
